# Tagihan Siswa - fixing bug laporan, pembayaran ujian, jenis ujian, revisi database
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data Pembayaran 1")

# Kelas 10 exam fee section (rows 14-18)
$ws.Range("B14").Value = 300000
$ws.Range("B15").Value = "0,-"
$ws.Range("B16").Value = 300000
$ws.Range("B17").Value = "0,-"
$ws.Range("B18").Value = 760000

# Kelas 11 exam fee section (rows 23-27)
$ws.Range("B23").Value = 300000
$ws.Range("B24").Value = 150000
$ws.Range("B25").Value = 300000
$ws.Range("B26").Value = 150000
$ws.Range("B27").Value = 1140000

# Totals section (rows 32-34)
$ws.Range("B32").Value = 760000
$ws.Range("B33").Value = 1140000
$ws.Range("B34").Value = 2100000
